# Auto-generated edit script: updates LeveProfit/Price columns (H-N) across
# multiple sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR) per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 164337.05
$ws.Cells.Item(15, 9).Value = 164337.05
$ws.Cells.Item(15, 11).Value = 493011.15
$ws.Cells.Item(15, 13).Value = -492842.15
$ws.Cells.Item(17, 8).Value = 1999
$ws.Cells.Item(17, 10).Value = 2176.5557
$ws.Cells.Item(17, 12).Value = 6529.6671
$ws.Cells.Item(17, 14).Value = -6865.6671
$ws.Cells.Item(41, 8).Value = 180.22728
$ws.Cells.Item(41, 9).Value = 127.333336
$ws.Cells.Item(41, 11).Value = 127.333336
$ws.Cells.Item(41, 13).Value = 312.666664
$ws.Cells.Item(61, 8).Value = 1223
$ws.Cells.Item(61, 9).Value = 705
$ws.Cells.Item(61, 11).Value = 2115
$ws.Cells.Item(61, 13).Value = -1943
$ws.Cells.Item(70, 8).Value = 6347.4
$ws.Cells.Item(70, 9).Value = 749.5
$ws.Cells.Item(70, 10).Value = 7746.875
$ws.Cells.Item(70, 11).Value = 2248.5
$ws.Cells.Item(70, 12).Value = 23240.625
$ws.Cells.Item(70, 13).Value = -1978.5
$ws.Cells.Item(70, 14).Value = -23780.625
$ws.Cells.Item(73, 8).Value = 6347.4
$ws.Cells.Item(73, 9).Value = 749.5
$ws.Cells.Item(73, 10).Value = 7746.875
$ws.Cells.Item(73, 11).Value = 2248.5
$ws.Cells.Item(73, 12).Value = 23240.625
$ws.Cells.Item(73, 13).Value = -1312.5
$ws.Cells.Item(73, 14).Value = -25112.625
$ws.Cells.Item(80, 8).Value = 1282.8334
$ws.Cells.Item(80, 9).Value = 803.2857
$ws.Cells.Item(80, 11).Value = 2409.8571
$ws.Cells.Item(80, 13).Value = -1411.8571
$ws.Cells.Item(83, 8).Value = 1282.8334
$ws.Cells.Item(83, 9).Value = 803.2857
$ws.Cells.Item(83, 11).Value = 7229.571300000001
$ws.Cells.Item(83, 13).Value = -2237.571300000001
$ws.Cells.Item(88, 8).Value = 2734.5
$ws.Cells.Item(88, 9).Value = 892
$ws.Cells.Item(88, 10).Value = 2876.2307
$ws.Cells.Item(88, 11).Value = 892
$ws.Cells.Item(88, 12).Value = 2876.2307
$ws.Cells.Item(88, 13).Value = -486
$ws.Cells.Item(88, 14).Value = -3688.2307
$ws.Cells.Item(91, 8).Value = 2734.5
$ws.Cells.Item(91, 9).Value = 892
$ws.Cells.Item(91, 10).Value = 2876.2307
$ws.Cells.Item(91, 11).Value = 892
$ws.Cells.Item(91, 12).Value = 2876.2307
$ws.Cells.Item(91, 13).Value = 512
$ws.Cells.Item(91, 14).Value = -5684.2307
$ws.Cells.Item(98, 8).Value = 3322.4285
$ws.Cells.Item(98, 9).Value = 3742.8333
$ws.Cells.Item(98, 10).Value = 800
$ws.Cells.Item(98, 11).Value = 3742.8333
$ws.Cells.Item(98, 12).Value = 800
$ws.Cells.Item(98, 13).Value = -2244.8333
$ws.Cells.Item(98, 14).Value = -3796
$ws.Cells.Item(122, 8).Value = 3322.4285
$ws.Cells.Item(122, 9).Value = 3742.8333
$ws.Cells.Item(122, 10).Value = 800
$ws.Cells.Item(122, 11).Value = 11228.4999
$ws.Cells.Item(122, 12).Value = 2400
$ws.Cells.Item(122, 13).Value = -8778.499899999999
$ws.Cells.Item(122, 14).Value = -7300
$ws.Cells.Item(132, 8).Value = 282735.06
$ws.Cells.Item(132, 9).Value = 306373.75
$ws.Cells.Item(132, 11).Value = 919121.25
$ws.Cells.Item(132, 13).Value = -916591.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 5472.0454
$ws.Cells.Item(63, 10).Value = 5849.1333
$ws.Cells.Item(63, 12).Value = 5849.1333
$ws.Cells.Item(63, 14).Value = -7221.1333
$ws.Cells.Item(66, 8).Value = 5472.0454
$ws.Cells.Item(66, 10).Value = 5849.1333
$ws.Cells.Item(66, 12).Value = 29245.6665
$ws.Cells.Item(66, 14).Value = -36109.66650000001
$ws.Cells.Item(74, 8).Value = 4082.1738
$ws.Cells.Item(74, 9).Value = 5882.75
$ws.Cells.Item(74, 11).Value = 5882.75
$ws.Cells.Item(74, 13).Value = -5008.75
$ws.Cells.Item(77, 8).Value = 4082.1738
$ws.Cells.Item(77, 9).Value = 5882.75
$ws.Cells.Item(77, 11).Value = 29413.75
$ws.Cells.Item(77, 13).Value = -25045.75
$ws.Cells.Item(132, 8).Value = 740683.9399999999
$ws.Cells.Item(132, 9).Value = 786592.8
$ws.Cells.Item(132, 10).Value = 235686.25
$ws.Cells.Item(132, 11).Value = 2359778.4
$ws.Cells.Item(132, 12).Value = 707058.75
$ws.Cells.Item(132, 13).Value = -2357248.4
$ws.Cells.Item(132, 14).Value = -712118.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 7577.9517
$ws.Cells.Item(99, 9).Value = 7272.116
$ws.Cells.Item(99, 11).Value = 7272.116
$ws.Cells.Item(99, 13).Value = -5774.116
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 11234.875
$ws.Cells.Item(62, 9).Value = 10999.167
$ws.Cells.Item(62, 10).Value = 11942
$ws.Cells.Item(62, 11).Value = 10999.167
$ws.Cells.Item(62, 12).Value = 11942
$ws.Cells.Item(62, 13).Value = -10375.167
$ws.Cells.Item(62, 14).Value = -13190
$ws.Cells.Item(65, 8).Value = 11234.875
$ws.Cells.Item(65, 9).Value = 10999.167
$ws.Cells.Item(65, 10).Value = 11942
$ws.Cells.Item(65, 11).Value = 54995.835
$ws.Cells.Item(65, 12).Value = 59710
$ws.Cells.Item(65, 13).Value = -51875.835
$ws.Cells.Item(65, 14).Value = -65950
$ws.Cells.Item(132, 8).Value = 18138.166
$ws.Cells.Item(132, 9).Value = 7002.8887
$ws.Cells.Item(132, 10).Value = 29273.445
$ws.Cells.Item(132, 11).Value = 21008.6661
$ws.Cells.Item(132, 12).Value = 87820.33499999999
$ws.Cells.Item(132, 13).Value = -18478.6661
$ws.Cells.Item(132, 14).Value = -92880.33499999999
$ws.Cells.Item(134, 8).Value = 66676704
$ws.Cells.Item(134, 10).Value = 18568.8
$ws.Cells.Item(134, 12).Value = 55706.39999999999
$ws.Cells.Item(134, 14).Value = -60776.39999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 6428.9473
$ws.Cells.Item(80, 9).Value = 4500.3335
$ws.Cells.Item(80, 10).Value = 8164.7
$ws.Cells.Item(80, 11).Value = 4500.3335
$ws.Cells.Item(80, 12).Value = 8164.7
$ws.Cells.Item(80, 13).Value = -3502.3335
$ws.Cells.Item(80, 14).Value = -10160.7
$ws.Cells.Item(83, 8).Value = 6428.9473
$ws.Cells.Item(83, 9).Value = 4500.3335
$ws.Cells.Item(83, 10).Value = 8164.7
$ws.Cells.Item(83, 11).Value = 22501.6675
$ws.Cells.Item(83, 12).Value = 40823.5
$ws.Cells.Item(83, 13).Value = -17509.6675
$ws.Cells.Item(83, 14).Value = -50807.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 5500
$ws.Cells.Item(14, 10).Value = 5500
$ws.Cells.Item(14, 12).Value = 5500
$ws.Cells.Item(14, 14).Value = -5844
$ws.Cells.Item(46, 8).Value = 35715090
$ws.Cells.Item(46, 9).Value = 649.5714
$ws.Cells.Item(46, 10).Value = 71429530
$ws.Cells.Item(46, 11).Value = 649.5714
$ws.Cells.Item(46, 12).Value = 71429530
$ws.Cells.Item(46, 13).Value = -461.5714
$ws.Cells.Item(46, 14).Value = -71429906
$ws.Cells.Item(55, 8).Value = 4589.0835
$ws.Cells.Item(55, 9).Value = 2509.9285
$ws.Cells.Item(55, 11).Value = 2509.9285
$ws.Cells.Item(55, 13).Value = -2336.9285
$ws.Cells.Item(68, 8).Value = 4062.125
$ws.Cells.Item(68, 9).Value = 3299.8572
$ws.Cells.Item(68, 10).Value = 4655
$ws.Cells.Item(68, 11).Value = 3299.8572
$ws.Cells.Item(68, 12).Value = 4655
$ws.Cells.Item(68, 13).Value = -2550.8572
$ws.Cells.Item(68, 14).Value = -6153
$ws.Cells.Item(71, 8).Value = 4062.125
$ws.Cells.Item(71, 9).Value = 3299.8572
$ws.Cells.Item(71, 10).Value = 4655
$ws.Cells.Item(71, 11).Value = 16499.286
$ws.Cells.Item(71, 12).Value = 23275
$ws.Cells.Item(71, 13).Value = -12755.286
$ws.Cells.Item(71, 14).Value = -30763
$ws.Cells.Item(136, 8).Value = 35722932
$ws.Cells.Item(136, 9).Value = 71436770
$ws.Cells.Item(136, 10).Value = 9096.429
$ws.Cells.Item(136, 11).Value = 214310310
$ws.Cells.Item(136, 12).Value = 27289.287
$ws.Cells.Item(136, 13).Value = -214307760
$ws.Cells.Item(136, 14).Value = -32389.287
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 7000000
$ws.Cells.Item(12, 10).Value = 7000000
$ws.Cells.Item(12, 12).Value = 7000000
$ws.Cells.Item(12, 14).Value = -7000284
